$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = @'
Pipeline(steps=[('scaler', None), ('selector', None),
                ('model',
                 MLPClassifier(activation='tanh', alpha=1e-05,
                               hidden_layer_sizes=(5, 10, 5),
                               learning_rate_init=0.0001, max_iter=1000,
                               random_state=42))])
'@
$ws.Range("B2").Value = 0.758095238095238
$ws.Range("C2").Value = @'
{'selector': None, 'scaler': None, 'model__solver': 'adam', 'model__learning_rate_init': 0.0001, 'model__hidden_layer_sizes': (5, 10, 5), 'model__alpha': 1e-05, 'model__activation': 'tanh'}
'@
$ws.Range("D2").Value = 0.7439022533136688
$ws.Range("E2").Value = 0.5466748640248641
$ws.Range("F2").Value = 0.7692307692307693
$ws.Range("G2").Value = 0.727466508871099
$ws.Range("H2").Value = 0.5356468253968254
$ws.Range("I2").Value = 0.6521739130434783
$ws.Range("J2").Value = 0.810936170212766
$ws.Range("K2").Value = 0.6096666666666666
$ws.Range("L2").Value = 0.9375
$ws.Range("M2").Value = @'
[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]
'@
$ws.Range("N2").Value = @'
[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 0]
'@

$ws.Range("A3").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()), ('selector', None),
                ('model',
                 MLPClassifier(hidden_layer_sizes=(10,),
                               learning_rate_init=1e-05, max_iter=1000,
                               random_state=42))])
'@
$ws.Range("B3").Value = 0.7538095238095237
$ws.Range("C3").Value = @'
{'selector': None, 'scaler': MinMaxScaler(), 'model__solver': 'adam', 'model__learning_rate_init': 1e-05, 'model__hidden_layer_sizes': (10,), 'model__alpha': 0.0001, 'model__activation': 'relu'}
'@
$ws.Range("D3").Value = 0.6999924646293084
$ws.Range("E3").Value = 0.5185700022200022
$ws.Range("F3").Value = 0.8
$ws.Range("G3").Value = 0.692439401571078
$ws.Range("H3").Value = 0.5339107142857142
$ws.Range("I3").Value = 0.6666666666666666
$ws.Range("J3").Value = 0.759404255319149
$ws.Range("K3").Value = 0.5561666666666667
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = @'
[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]
'@
$ws.Range("N3").Value = @'
[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]
'@

$ws.Range("A4").Value = @'
Pipeline(steps=[('scaler', StandardScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7faa9f72b580>),
                ('model',
                 MLPClassifier(alpha=1, hidden_layer_sizes=(10,),
                               learning_rate_init=0.01, max_iter=1000,
                               random_state=42, solver='lbfgs'))])
'@
$ws.Range("B4").Value = 0.7004545454545454
$ws.Range("C4").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7faa00146850>, 'scaler': StandardScaler(), 'model__solver': 'lbfgs', 'model__learning_rate_init': 0.01, 'model__hidden_layer_sizes': (10,), 'model__alpha': 1, 'model__activation': 'relu'}
'@
$ws.Range("D4").Value = 0.6711423221211917
$ws.Range("E4").Value = 0.4897046287046289
$ws.Range("F4").Value = 0.7096774193548387
$ws.Range("G4").Value = 0.6784388231901994
$ws.Range("H4").Value = 0.4795238095238094
$ws.Range("I4").Value = 0.9166666666666666
$ws.Range("J4").Value = 0.7231555555555554
$ws.Range("K4").Value = 0.5666
$ws.Range("L4").Value = 0.5789473684210527
$ws.Range("M4").Value = @'
[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]
'@
$ws.Range("N4").Value = @'
[0 1 1 0 0 1 0 1 1 0 0 0 0 1 1 1 1 1 1 0 0 0 1 0]
'@

$ws.Range("A5").Value = @'
Pipeline(steps=[('scaler', MinMaxScaler()), ('selector', None),
                ('model',
                 MLPClassifier(activation='tanh', alpha=0.01,
                               hidden_layer_sizes=(10, 10, 10),
                               learning_rate_init=0.0001, max_iter=1000,
                               random_state=42))])
'@
$ws.Range("B5").Value = 0.7699999999999999
$ws.Range("C5").Value = @'
{'selector': None, 'scaler': MinMaxScaler(), 'model__solver': 'adam', 'model__learning_rate_init': 0.0001, 'model__hidden_layer_sizes': (10, 10, 10), 'model__alpha': 0.01, 'model__activation': 'tanh'}
'@
$ws.Range("D5").Value = 0.7200633141986819
$ws.Range("E5").Value = 0.5592921911421911
$ws.Range("F5").Value = 0.7368421052631579
$ws.Range("G5").Value = 0.6923024139416202
$ws.Range("H5").Value = 0.5351384920634922
$ws.Range("I5").Value = 0.5833333333333334
$ws.Range("J5").Value = 0.7887959183673469
$ws.Range("K5").Value = 0.6288333333333332
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = @'
[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]
'@
$ws.Range("N5").Value = @'
[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]
'@

$ws.Range("A6").Value = @'
Pipeline(steps=[('scaler', StandardScaler()),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 MLPClassifier(activation='tanh', alpha=0.01,
                               hidden_layer_sizes=(5, 10, 5),
                               learning_rate_init=1, max_iter=1000,
                               random_state=42))])
'@
$ws.Range("B6").Value = 0.7592424242424242
$ws.Range("C6").Value = @'
{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': StandardScaler(), 'model__solver': 'adam', 'model__learning_rate_init': 1, 'model__hidden_layer_sizes': (5, 10, 5), 'model__alpha': 0.01, 'model__activation': 'tanh'}
'@
$ws.Range("D6").Value = 0.7389922420838458
$ws.Range("E6").Value = 0.5598278332778333
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.7238216224078028
$ws.Range("H6").Value = 0.5380253968253969
$ws.Range("I6").Value = 0.4117647058823529
$ws.Range("J6").Value = 0.7931153846153847
$ws.Range("K6").Value = 0.6273333333333334
$ws.Range("L6").Value = 0.6363636363636364
$ws.Range("M6").Value = @'
[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]
'@
$ws.Range("N6").Value = @'
[1 1 1 0 1 1 1 1 0 0 0 1 1 1 1 0 0 1 0 1 1 1 1 1]
'@
